# Auto-generated Excel COM-interop edit script
# Applies the numeric corrections described in the commit diff
# for the 'Siren_Profits' scheduled-runner update across all sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 265.66666
$ws.Range("I11").Value = 265.66666
$ws.Range("K11").Value = 265.66666
$ws.Range("M11").Value = -125.66666
# Row 33
$ws.Range("H33").Value = 459.79166
$ws.Range("I33").Value = 261.8421
$ws.Range("J33").Value = 1212
$ws.Range("K33").Value = 261.8421
$ws.Range("L33").Value = 1212
$ws.Range("M33").Value = -32.84210000000002
$ws.Range("N33").Value = -1670
# Row 76
$ws.Range("H76").Value = 5038.5
$ws.Range("I76").Value = 3828.4285
$ws.Range("K76").Value = 3828.4285
$ws.Range("M76").Value = -3513.4285
# Row 79
$ws.Range("H79").Value = 5038.5
$ws.Range("I79").Value = 3828.4285
$ws.Range("K79").Value = 3828.4285
$ws.Range("M79").Value = -2736.4285
# Row 86
$ws.Range("H86").Value = 47643612
$ws.Range("I86").Value = 2785.4
$ws.Range("K86").Value = 2785.4
$ws.Range("M86").Value = -1662.4
# Row 89
$ws.Range("H89").Value = 47643612
$ws.Range("I89").Value = 2785.4
$ws.Range("K89").Value = 13927
$ws.Range("M89").Value = -8311
# Row 132
$ws.Range("H132").Value = 3875.195
$ws.Range("I132").Value = 3664.4055
$ws.Range("K132").Value = 10993.2165
$ws.Range("M132").Value = -8463.216499999999
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()
# Row 138
$ws.Range("H138").Value = 2297.7144
$ws.Range("J138").Value = 4117.8335
$ws.Range("L138").Value = 12353.5005
$ws.Range("N138").Value = -22633.5005
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3085.318
$ws.Range("I2").Value = 3060.389
$ws.Range("K2").Value = 3060.389
$ws.Range("M2").Value = -2947.389
# Row 32
$ws.Range("H32").Value = 7079.4116
$ws.Range("I32").Value = 7021
$ws.Range("K32").Value = 7021
$ws.Range("M32").Value = -6734
# Row 45
$ws.Range("H45").Value = 65049.21
$ws.Range("I45").Value = 129451.56
$ws.Range("K45").Value = 129451.56
$ws.Range("M45").Value = -129074.56
# Row 61
$ws.Range("H61").Value = 5773.353
$ws.Range("I61").Value = 6808.5884
$ws.Range("J61").Value = 3702.8823
$ws.Range("K61").Value = 6808.5884
$ws.Range("L61").Value = 3702.8823
$ws.Range("M61").Value = -6596.5884
$ws.Range("N61").Value = -4126.8823
# Row 96
$ws.Range("H96").Value = 35000
$ws.Range("J96").Value = 35000
$ws.Range("L96").Value = 35000
$ws.Range("N96").Value = -40492
# Row 102
$ws.Range("H102").Value = 7941.212
$ws.Range("I102").Value = 9802.32
$ws.Range("J102").Value = 2125.25
$ws.Range("K102").Value = 9802.32
$ws.Range("L102").Value = 2125.25
$ws.Range("M102").Value = -8180.32
$ws.Range("N102").Value = -5369.25
# Row 116
$ws.Range("H116").Value = 3085.318
$ws.Range("I116").Value = 3060.389
$ws.Range("K116").Value = 3060.389
$ws.Range("M116").Value = -766.3890000000001
# Row 136
$ws.Range("H136").Value = 5773.353
$ws.Range("I136").Value = 6808.5884
$ws.Range("J136").Value = 3702.8823
$ws.Range("K136").Value = 20425.7652
$ws.Range("L136").Value = 11108.6469
$ws.Range("M136").Value = -17875.7652
$ws.Range("N136").Value = -16208.6469
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3085.318
$ws.Range("I3").Value = 3060.389
$ws.Range("K3").Value = 3060.389
$ws.Range("M3").Value = -2946.389
# Row 86
$ws.Range("H86").Value = 4996.0415
$ws.Range("I86").Value = 6912
$ws.Range("K86").Value = 6912
$ws.Range("M86").Value = -5789
# Row 89
$ws.Range("H89").Value = 4996.0415
$ws.Range("I89").Value = 6912
$ws.Range("K89").Value = 34560
$ws.Range("M89").Value = -28944
# Row 134
$ws.Range("H134").Value = 8132.45
$ws.Range("I134").Value = 9604.625
$ws.Range("J134").Value = 2243.75
$ws.Range("K134").Value = 28813.875
$ws.Range("L134").Value = 6731.25
$ws.Range("M134").Value = -26278.875
$ws.Range("N134").Value = -11801.25
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8773.591
$ws.Range("I31").Value = 15449.444
$ws.Range("J31").Value = 4151.846
$ws.Range("K31").Value = 15449.444
$ws.Range("L31").Value = 4151.846
$ws.Range("M31").Value = -15154.444
$ws.Range("N31").Value = -4741.846
# Row 34
$ws.Range("H34").Value = 8773.591
$ws.Range("I34").Value = 15449.444
$ws.Range("J34").Value = 4151.846
$ws.Range("K34").Value = 15449.444
$ws.Range("L34").Value = 4151.846
$ws.Range("M34").Value = -15247.444
$ws.Range("N34").Value = -4555.846
# Row 58
$ws.Range("H58").Value = 2507.4473
$ws.Range("I58").Value = 2410.2964
$ws.Range("J58").Value = 2745.9092
$ws.Range("K58").Value = 2410.2964
$ws.Range("L58").Value = 2745.9092
$ws.Range("M58").Value = -2207.2964
$ws.Range("N58").Value = -3151.9092
# Row 69
$ws.Range("H69").Value = 4750
$ws.Range("I69").Value = 4750
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 4750
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -4001
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 4750
$ws.Range("I72").Value = 4750
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 14250
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -10506
$ws.Range("N72").ClearContents()
# Row 134
$ws.Range("H134").Value = 7281.826
$ws.Range("I134").Value = 9311.5625
$ws.Range("K134").Value = 27934.6875
$ws.Range("M134").Value = -25399.6875
# Row 136
$ws.Range("H136").Value = 2507.4473
$ws.Range("I136").Value = 2410.2964
$ws.Range("J136").Value = 2745.9092
$ws.Range("K136").Value = 7230.889200000001
$ws.Range("L136").Value = 8237.7276
$ws.Range("M136").Value = -4680.889200000001
$ws.Range("N136").Value = -13337.7276
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 587.3125
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
# Row 55
$ws.Range("H55").Value = 5701.8887
$ws.Range("J55").Value = 6866.143
$ws.Range("L55").Value = 20598.429
$ws.Range("N55").Value = -20952.429
# Row 58
$ws.Range("H58").Value = 3041.8462
$ws.Range("I58").Value = 522
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 1566
$ws.Range("L58").Value = 10500
$ws.Range("M58").Value = -1438
$ws.Range("N58").Value = -10756
# Row 61
$ws.Range("H61").Value = 366.85715
$ws.Range("I61").Value = 223.33333
$ws.Range("K61").Value = 669.99999
$ws.Range("M61").Value = -454.99999
# Row 81
$ws.Range("H81").Value = 7322.1113
$ws.Range("I81").Value = 2731
$ws.Range("K81").Value = 8193
$ws.Range("M81").Value = -7070
# Row 84
$ws.Range("H84").Value = 7322.1113
$ws.Range("I84").Value = 2731
$ws.Range("K84").Value = 24579
$ws.Range("M84").Value = -18963
# Row 105
$ws.Range("H105").Value = 9951.829
$ws.Range("I105").Value = 8026
$ws.Range("J105").Value = 9999.975
$ws.Range("K105").Value = 24078
$ws.Range("L105").Value = 29999.925
$ws.Range("M105").Value = -21457
$ws.Range("N105").Value = -35241.925
# Row 134
$ws.Range("H134").Value = 3222
$ws.Range("I134").Value = 2968
$ws.Range("K134").Value = 8904
$ws.Range("M134").Value = -3834
# Row 138
$ws.Range("H138").Value = 913.625
$ws.Range("I138").Value = 925.1429000000001
$ws.Range("K138").Value = 2775.4287
$ws.Range("M138").Value = 2364.5713
# Row 139
$ws.Range("H139").Value = 1113434.8
$ws.Range("I139").Value = 2000982.5
$ws.Range("K139").Value = 6002947.5
$ws.Range("M139").Value = -5997807.5
# Row 140
$ws.Range("H140").Value = 11069.111
$ws.Range("I140").Value = 11440.823
$ws.Range("K140").Value = 34322.469
$ws.Range("M140").Value = -29142.469
# Row 141
$ws.Range("H141").Value = 2544.4443
$ws.Range("I141").Value = 2425
$ws.Range("K141").Value = 7275
$ws.Range("M141").Value = -2095
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6195.8237
$ws.Range("I70").Value = 5496.5
$ws.Range("K70").Value = 5496.5
$ws.Range("M70").Value = -5226.5
# Row 73
$ws.Range("H73").Value = 6195.8237
$ws.Range("I73").Value = 5496.5
$ws.Range("K73").Value = 5496.5
$ws.Range("M73").Value = -4560.5
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3430.8667
$ws.Range("I16").Value = 3247.3572
$ws.Range("J16").Value = 6000
$ws.Range("K16").Value = 3247.3572
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = -3077.3572
$ws.Range("N16").Value = -6340
# Row 93
$ws.Range("H93").Value = 3486.1667
$ws.Range("I93").Value = 4029.875
$ws.Range("K93").Value = 4029.875
$ws.Range("M93").Value = -2781.875
# Row 122
$ws.Range("H122").Value = 4583.525
$ws.Range("I122").Value = 4314.8066
$ws.Range("J122").Value = 5509.1113
$ws.Range("K122").Value = 12944.4198
$ws.Range("L122").Value = 16527.3339
$ws.Range("M122").Value = -10494.4198
$ws.Range("N122").Value = -21427.3339
# Row 139
$ws.Range("H139").Value = 68749.5
$ws.Range("J139").Value = 68749.5
$ws.Range("L139").Value = 68749.5
$ws.Range("N139").Value = -79029.5
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 25134
$ws.Range("I107").Value = 2475.5
$ws.Range("J107").Value = 100662.336
$ws.Range("K107").Value = 7426.5
$ws.Range("L107").Value = 301987.008
$ws.Range("M107").Value = -5506.5
$ws.Range("N107").Value = -305827.008
# Row 122
$ws.Range("H122").Value = 4364.5713
$ws.Range("I122").Value = 2007.7222
$ws.Range("K122").Value = 6023.1666
$ws.Range("M122").Value = -3573.1666
# Row 132
$ws.Range("H132").Value = 6402.6772
$ws.Range("I132").Value = 7382.022
$ws.Range("K132").Value = 22146.066
$ws.Range("M132").Value = -19616.066
